# Fruta / hortaliza, semanal
# Insert a new weekly data row in chronological position (new row 5),
# pushing the existing rows 5-13 down to rows 6-14, and append one more
# new row at the end (row 15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 5; this shifts old rows 5-13
# down to 6-14 and copies formatting (e.g. date style on column D) from
# the adjacent row, matching Excel's native "Insert Row" behavior.
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the new weekly record.
$ws.Range("A5").Value2 = 10
$ws.Range("B5").Value2 = "Vega Modelo de Temuco"
$ws.Range("C5").Value2 = "La Araucanía"
$ws.Range("D5").Value2 = 44424
$ws.Range("E5").Value2 = 9
$ws.Range("F5").Value2 = "Fruta"
$ws.Range("G5").Value2 = 100108
$ws.Range("H5").Value2 = "Tropicales y subtropicales"
$ws.Range("I5").Value2 = 100108003
$ws.Range("J5").Value2 = "Maracuyá"
$ws.Range("K5").Value2 = "Sin especificar"
$ws.Range("L5").Value2 = "Primera"
$ws.Range("M5").Value2 = 15
$ws.Range("N5").Value2 = 35000
$ws.Range("O5").Value2 = 35000
$ws.Range("P5").Value2 = 35000
$ws.Range("Q5").Value2 = "$/caja 18 kilos"
$ws.Range("R5").Value2 = "Región de Arica y Parinacota"
$ws.Range("S5").Value2 = 1944
$ws.Range("T5").Value2 = 18

# Append a brand-new row (15) with another weekly record.
$ws.Range("A15").Value2 = 10
$ws.Range("B15").Value2 = "Vega Modelo de Temuco"
$ws.Range("C15").Value2 = "La Araucanía"
$ws.Range("D15").Value2 = 44418
$ws.Range("E15").Value2 = 9
$ws.Range("F15").Value2 = "Fruta"
$ws.Range("G15").Value2 = 100108
$ws.Range("H15").Value2 = "Tropicales y subtropicales"
$ws.Range("I15").Value2 = 100108003
$ws.Range("J15").Value2 = "Maracuyá"
$ws.Range("K15").Value2 = "Sin especificar"
$ws.Range("L15").Value2 = "Primera"
$ws.Range("M15").Value2 = 30
$ws.Range("N15").Value2 = 35000
$ws.Range("O15").Value2 = 35000
$ws.Range("P15").Value2 = 35000
$ws.Range("Q15").Value2 = "$/caja 18 kilos"
$ws.Range("R15").Value2 = "Región de Arica y Parinacota"
$ws.Range("S15").Value2 = 1944
$ws.Range("T15").Value2 = 18

# Make sure the new date cell uses the same date style as the rest of
# column D (row-insert should already have copied it, but set it
# explicitly to be safe).
$ws.Range("D15").NumberFormat = $ws.Range("D14").NumberFormat
